$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric but must stay stored as text
# (matches the source data which keeps these as literal strings).
$textForced = [ordered]@{
    'D5' = '214.71'
    'D8' = '0.0628'
    'D9' = '0.250'
    'D16' = '64.82'
    'D19' = '213.64'
    'D22' = '4.35'
    'D24' = '9.05'
    'D25' = '148.23'
    'D27' = '7.37'
    'D32' = '0.772'
    'D33' = '3.34'
    'D36' = '1.57'
    'D41' = '0.801'
    'D43' = '64.98'
    'D44' = '5.34'
    'D46' = '0.880'
    'D47' = '89.97'
    'D48' = '1.64'
    'D50' = '0.102'
}
foreach ($ref in $textForced.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textForced[$ref]
    $cell.Style = "Normal"
}

# Remaining plain text / percentage / label / URL updates.
$plainValues = [ordered]@{
    'D2' = '27.006.40'
    'E2' = '  -0.67%  '
    'D3' = '1.620.68'
    'E3' = '  -1.28%  '
    'E4' = '  -0.15%  '
    'E5' = '  -1.13%  '
    'E6' = '  -1.39%  '
    'B8' = 'Dogecoin'
    'C8' = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
    'E8' = '  +0.06%  '
    'B9' = 'Cardano'
    'C9' = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
    'E9' = '  -1.47%  '
    'E10' = '  +0.60%  '
    'E11' = '  -0.05%  '
    'D12' = '1.847.80'
    'E12' = '  -1.29%  '
    'D13' = '1.621.43'
    'E13' = '  -2.30%  '
    'E14' = '  +0.01%  '
    'E15' = '  -0.48%  '
    'E16' = '  -3.44%  '
    'D17' = '26.984.55'
    'E17' = '  -0.76%  '
    'E18' = '  +0.74%  '
    'E19' = '  -2.31%  '
    'E20' = '  -0.06%  '
    'E21' = '  -1.51%  '
    'E22' = '  -1.35%  '
    'E23' = '  -5.62%  '
    'E24' = '  -0.57%  '
    'E25' = '  +0.47%  '
    'E26' = '  -0.19%  '
    'E27' = '  -0.96%  '
    'E28' = '  -2.08%  '
    'E30' = '  +0.96%  '
    'E31' = '  -0.98%  '
    'B32' = 'ImmutableX'
    'C32' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'E32' = '  +40.39%  '
    'B33' = 'Filecoin'
    'C33' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'E33' = '  -1.16%  '
    'E34' = '  -0.08%  '
    'D35' = '1.348.69'
    'E35' = '  +3.43%  '
    'E36' = '  -0.64%  '
    'E37' = '  -0.51%  '
    'E38' = '  +0.57%  '
    'E39' = '  -1.37%  '
    'E40' = '  -0.14%  '
    'E41' = '  -1.14%  '
    'E42' = '  -0.20%  '
    'E43' = '  +5.22%  '
    'E44' = '  +0.29%  '
    'D45' = '1.758.89'
    'E45' = '  -1.34%  '
    'E46' = '  +31.67%  '
    'E47' = '  -1.95%  '
    'E48' = '  +2.02%  '
    'E49' = '  -1.63%  '
    'E50' = '  +5.69%  '
    'E51' = '  +0.41%  '
}
foreach ($ref in $plainValues.Keys) {
    $ws.Range($ref).Value = $plainValues[$ref]
}
